$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 1.02
    "C2" = 1.052042066530208
    "D2" = 1.055268628335639
    "E2" = 1.048900163478012
    "F2" = 1.063662236074625
    "I2" = 1.045481499151993
    "J2" = 1.057066755259954
    "K2" = 1.058009552368777
    "L2" = 1.051658720252838
    "M2" = 1.066380269394519
    "N2" = 1.005712725503983
    "B3" = 1.02
    "C3" = 1.053202855792059
    "D3" = 1.056180841177848
    "E3" = 1.049895937688397
    "F3" = 1.06477227772745
    "I3" = 1.045816074696081
    "J3" = 1.057876841161164
    "K3" = 1.058734920830761
    "L3" = 1.052466157880417
    "M3" = 1.06730462547779
    "B4" = 1.02
    "C4" = 1.053953946110288
    "D4" = 1.056771069751885
    "E4" = 1.050540534742435
    "F4" = 1.065490930156557
    "I4" = 1.046031356457702
    "J4" = 1.058400449160515
    "K4" = 1.059203617025321
    "L4" = 1.052988275643255
    "M4" = 1.06790253387806
    "B5" = 1.02
    "C5" = 1.054269701245172
    "D5" = 1.057019193805947
    "E5" = 1.050811586922102
    "F5" = 1.065793143309404
    "I5" = 1.046121571142185
    "J5" = 1.058620437496409
    "K5" = 1.059400497801696
    "L5" = 1.053207690850608
    "M5" = 1.068153844224088
    "B6" = 1.02
    "C6" = 1.054322717763011
    "D6" = 1.057060854458923
    "E6" = 1.050857101511107
    "F6" = 1.06584389161488
    "I6" = 1.046136701596776
    "J6" = 1.058657366515479
    "K6" = 1.059433545628235
    "L6" = 1.05324452676515
    "M6" = 1.068196037395194
    "B7" = 1.02
    "C7" = 1.053958165257223
    "D7" = 1.056774385229057
    "E7" = 1.0505441563048
    "F7" = 1.065494967982785
    "I7" = 1.046032563049137
    "J7" = 1.058403389191132
    "K7" = 1.059206248380725
    "L7" = 1.052991207806143
    "M7" = 1.067905892096575
    "B8" = 1.02
    "C8" = 1.052434365052744
    "D8" = 1.055576922084551
    "E8" = 1.049236635403579
    "F8" = 1.064037301147378
    "I8" = 1.045594821353602
    "J8" = 1.057340646214249
    "K8" = 1.058254832036046
    "L8" = 1.051931669894106
    "M8" = 1.066692703930327
    "B9" = 1.02
    "C9" = 1.049749049066064
    "D9" = 1.053466571125703
    "E9" = 1.046934631300531
    "F9" = 1.061471607575361
    "I9" = 1.044814183154003
    "J9" = 1.055463559655931
    "K9" = 1.056573204029831
    "L9" = 1.050061946698135
    "M9" = 1.0645532727909
    "B10" = 1.02
    "C10" = 1.047958644433357
    "D10" = 1.052059478447176
    "E10" = 1.045401297752923
    "F10" = 1.059763067540757
    "I10" = 1.044287506706357
    "J10" = 1.054209183296155
    "K10" = 1.055448662871111
    "L10" = 1.048813639485408
    "M10" = 1.063125857977395
    "B11" = 1.02
    "C11" = 1.047183317225147
    "D11" = 1.051450142685977
    "E11" = 1.044737659174127
    "F11" = 1.059023699520276
    "I11" = 1.0440579645218
    "J11" = 1.053665309286691
    "K11" = 1.054960899367457
    "L11" = 1.048272670169406
    "M11" = 1.062507497647956
    "B12" = 1.02
    "C12" = 1.046895314133787
    "D12" = 1.051223799529002
    "E12" = 1.044491199536419
    "F12" = 1.05874913085244
    "I12" = 1.043972478499233
    "J12" = 1.053463181266122
    "K12" = 1.054779597089528
    "L12" = 1.048071662785501
    "M12" = 1.062277768311135
    "B13" = 1.02
    "C13" = 1.04695709234366
    "D13" = 1.051272351305911
    "E13" = 1.04454406392503
    "F13" = 1.058808023809906
    "I13" = 1.043990825675208
    "J13" = 1.053506543356444
    "K13" = 1.054818492717484
    "L13" = 1.04811478260615
    "M13" = 1.06232704798022
    "B14" = 1.02
    "C14" = 1.047159511054152
    "D14" = 1.051431433270365
    "E14" = 1.044717285823032
    "F14" = 1.05900100225588
    "I14" = 1.044050902788451
    "J14" = 1.053648603549078
    "K14" = 1.054945915421627
    "L14" = 1.048256056212459
    "M14" = 1.062488509025823
    "B15" = 1.02
    "C15" = 1.047284226261454
    "D15" = 1.05152944775067
    "E15" = 1.044824019535208
    "F15" = 1.059119911307668
    "I15" = 1.044087888611214
    "J15" = 1.053736117057876
    "K15" = 1.055024408149639
    "L15" = 1.048343090609587
    "M15" = 1.062587984816467
    "B16" = 1.02
    "C16" = 1.048010098739864
    "D16" = 1.052099916853851
    "E16" = 1.045445347639079
    "F16" = 1.059812146191881
    "I16" = 1.044302709292409
    "J16" = 1.054245263176581
    "K16" = 1.055481016600854
    "L16" = 1.048849532432126
    "M16" = 1.063166890553526
    "B17" = 1.02
    "C17" = 1.048465399518386
    "D17" = 1.052457741944658
    "E17" = 1.045835171611753
    "F17" = 1.060246484421029
    "I17" = 1.044437062083937
    "J17" = 1.05456444361427
    "K17" = 1.055767212492498
    "L17" = 1.049167090712891
    "M17" = 1.06352994731451
    "B18" = 1.02
    "C18" = 1.048730962087931
    "D18" = 1.052666449953691
    "E18" = 1.046062578800715
    "F18" = 1.060499869169618
    "I18" = 1.044515284262995
    "J18" = 1.054750546667459
    "K18" = 1.055934065689893
    "L18" = 1.049352274279382
    "M18" = 1.063741685181042
    "B19" = 1.02
    "C19" = 1.048821510975428
    "D19" = 1.052737613102604
    "E19" = 1.04614012378888
    "F19" = 1.060586274070177
    "I19" = 1.044541931668639
    "J19" = 1.054813991233132
    "K19" = 1.055990944761047
    "L19" = 1.04941540984694
    "M19" = 1.063813877762384
    "B20" = 1.02
    "C20" = 1.04841655075445
    "D20" = 1.052419351251151
    "E20" = 1.04579334414519
    "F20" = 1.060199879637562
    "I20" = 1.044422662149165
    "J20" = 1.054530205767276
    "K20" = 1.055736514662843
    "L20" = 1.049133024146375
    "M20" = 1.063490997557358
    "B21" = 1.02
    "C21" = 1.047099904133601
    "D21" = 1.051384587840984
    "E21" = 1.044666275045923
    "F21" = 1.058944173110197
    "I21" = 1.044033217757394
    "J21" = 1.053606773377672
    "K21" = 1.054908396071458
    "L21" = 1.04821445651883
    "M21" = 1.06244096395844
    "B22" = 1.02
    "C22" = 1.046272005716384
    "D22" = 1.050733940333858
    "E22" = 1.043957903672615
    "F22" = 1.058155039381531
    "I22" = 1.043787063419849
    "J22" = 1.053025543806108
    "K22" = 1.054387000441574
    "L22" = 1.04763652696841
    "M22" = 1.061780518328757
    "B23" = 1.02
    "C23" = 1.046710897508291
    "D23" = 1.051078865829261
    "E23" = 1.044333400114326
    "F23" = 1.058573338440847
    "I23" = 1.043917677363587
    "J23" = 1.053333724645841
    "K23" = 1.0546634709922
    "L23" = 1.047942935483397
    "M23" = 1.062130656767331
    "B24" = 1.02
    "C24" = 1.048438623428382
    "D24" = 1.052436698369246
    "E24" = 1.045812244085623
    "F24" = 1.060220938201198
    "I24" = 1.044429169303208
    "J24" = 1.054545676591053
    "K24" = 1.055750385937734
    "L24" = 1.049148417493436
    "M24" = 1.063508597359207
    "B25" = 1.02
    "C25" = 1.050443296726022
    "D25" = 1.054012180401446
    "E25" = 1.04752951770538
    "F25" = 1.062134559941985
    "I25" = 1.045017097837978
    "J25" = 1.055949355611566
    "K25" = 1.057008552853121
    "L25" = 1.050545635040737
    "M25" = 1.065106563847094
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
